$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 156
$ws.Cells.Item(156, 1).Value = 155
$ws.Cells.Item(156, 2).Value = "Aysén"
$ws.Cells.Item(156, 3).Value = "Indicadores de desempeño"
$ws.Cells.Item(156, 4).Value = "Están trabajando para hacer un plan de trabajo más o menos rápido, para ponerse al día con las metas. Primeros avances, la semana pasada Jefe Prevención salió a un recorrido por varios lugares, a fin de aumentar capacitados y mejorar procesos donde se está más débil. Sílice, está lento, pero que ya se tienen los lugares y tomando las medias."
$ws.Cells.Item(156, 5).Value = "Pendiente"
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = "16-09-2025"

# Row 157
$ws.Cells.Item(157, 1).Value = 156
$ws.Cells.Item(157, 2).Value = "Aysén"
$ws.Cells.Item(157, 3).Value = "Ejecución Presupuestaria"
$ws.Cells.Item(157, 4).Value = "Faltan un par de compras que se han ido retrasando por distintas situaciones.
Por ejemplo, unas chaquetas para prevención, pero hubo que bajar la licitación y modifica las bases.
El resto se encuentra en regla. Salió lo de las pausas que era una preocupación, pero que ya eso se subsanó (Julio-Agosto)"
$ws.Cells.Item(157, 5).Value = "Pendiente"
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = "16-09-2025"

# Row 158
$ws.Cells.Item(158, 1).Value = 157
$ws.Cells.Item(158, 2).Value = "Aysén"
$ws.Cells.Item(158, 3).Value = "Otros"
$ws.Cells.Item(158, 4).Value = "Correos desde Chile Chico respecto de los exámenes ocupacionales que hace el hospital allá. Ellos no han tenido novedades respecto de cómo poder facturar esas atenciones y se han mandado los correos a la Camila con copia a María Isabel
Licitación de evaluaciones ocupacionales y vigilancia (Se requiere Información)"
$ws.Cells.Item(158, 5).Value = "Pendiente"
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = "16-09-2025"

# Row 159
$ws.Cells.Item(159, 1).Value = 158
$ws.Cells.Item(159, 2).Value = "Aysén"
$ws.Cells.Item(159, 3).Value = "Infraestructura"
$ws.Cells.Item(159, 4).Value = "Aún está en proceso el cambio de inmueble. El proceso estuvo detenido un mes (Carolina Barrera estuvo fuera con licencia). Al propieatario se le han dado una serie de plazos que finalmente no se cumplieron por responsabilidad de ambos. (corredora ahí envió la documentación a último).
Seremi de Gobierno está interesada en el inmueble y ya cuenta con la autorización de DIPRES."
$ws.Cells.Item(159, 5).Value = "Pendiente"
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = "16-09-2025"

# Row 160
$ws.Cells.Item(160, 1).Value = 159
$ws.Cells.Item(160, 2).Value = "Aysén"
$ws.Cells.Item(160, 3).Value = "Temas de Personas"
$ws.Cells.Item(160, 4).Value = "Jimena Zárate presentó una DIEP (salud mental). Salió como común. 
Clima Laboral: Aunque los problemas son derivados de asuntos personales, se ha logrado disminuir los conflictos y mantener un buen clima laboral en el último tiempo."
$ws.Cells.Item(160, 5).Value = "Pendiente"
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = "16-09-2025"

# Row 161
$ws.Cells.Item(161, 1).Value = 160
$ws.Cells.Item(161, 2).Value = "O'Higgins"
$ws.Cells.Item(161, 3).Value = "Indicadores de desempeño"
$ws.Cells.Item(161, 4).Value = "Capacitados: Meta 4400. a la fecha aprox. 2300
Región indica que lleva un avance de 61%, lo que contrasta con el 51% según datos de planilla de preveción NC.
Se solicita envíen planilla para trabajar un una planilla o instrumento único para todas las regiones.
Problemas con tablets ya han sido superados en su totalidad."
$ws.Cells.Item(161, 5).Value = "Pendiente"
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = "16-09-2025"

# Row 162
$ws.Cells.Item(162, 1).Value = 161
$ws.Cells.Item(162, 2).Value = "O'Higgins"
$ws.Cells.Item(162, 3).Value = "Plan de SSPP"
$ws.Cells.Item(162, 4).Value = "Se está armando el Cosoc. Se han realido ya dos reuniones. Manual Ahumada está colaborando en eso.
SLEP Cachapoal. En estos momentos cuenta con 15 funcionarios administrativos (todos en funciones online) los que no han sido posible ubicar y contactar. Se buscará la alternativa por intermedio de la Seremía de Educación para tomar contacto con ellos.
SLEP cuenta con 79 colegio y una estimación de 12000 funcionario según información del SEREMI de educación."
$ws.Cells.Item(162, 5).Value = "Pendiente"
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = "16-09-2025"

# Row 163
$ws.Cells.Item(163, 1).Value = 162
$ws.Cells.Item(163, 2).Value = "O'Higgins"
$ws.Cells.Item(163, 3).Value = "Temas de Personas"
$ws.Cells.Item(163, 4).Value = "Sumario a funcionario Julio Urzua"
$ws.Cells.Item(163, 5).Value = "Pendiente"
$ws.Cells.Item(163, 6).Value = 0
$ws.Cells.Item(163, 7).Value = "16-09-2025"

Write-Output "Added rows 156-163"